$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(7)

# ---- Header row (B1:K1) ----
$ws.Cells.Item(1,2).Value = "company"
$ws.Cells.Item(1,3).Value = "name"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "property_category"
$ws.Cells.Item(1,6).Value = "category"
$ws.Cells.Item(1,7).Value = "date"
$ws.Cells.Item(1,8).Value = "legislator_name"
$ws.Cells.Item(1,9).Value = "legislator_id"
$ws.Cells.Item(1,10).Value = "source_file"
$ws.Cells.Item(1,11).Value = "index"

# copy header formatting (bold / border / centered) onto the newly-introduced header cells F1:K1
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:K1").PasteSpecial(-4122) | Out-Null

# ---- Column B (company) for data rows 2:5 ----
$ws.Cells.Item(2,2).Value = "南山人壽"
$ws.Cells.Item(3,2).Value = "南山人壽"
$ws.Cells.Item(4,2).Value = "南山人壽"
$ws.Cells.Item(5,2).Value = "南山人壽"

# ---- Column C (policy name) for data rows 2:5 ----
$ws.Cells.Item(2,3).Value = "新20年期缴費增值分紅终身壽險南山終身醫療保險"
$ws.Cells.Item(3,3).Value = "南山新年年春還本終身保險南山终身醫療保險"
$ws.Cells.Item(4,3).Value = "南山新年年春還本終身保險南山终身醫療保險"
$ws.Cells.Item(5,3).Value = "南山終身醫療保險"

# ---- Column D (owner) for data rows 2:5 ----
$ws.Cells.Item(2,4).Value = "饒月琴"
$ws.Cells.Item(3,4).Value = "饒月琴"
$ws.Cells.Item(4,4).Value = "饒月琴"
$ws.Cells.Item(5,4).Value = "饒月琴"

# ---- Column E (property_category) for data rows 2:5 ----
$ws.Cells.Item(2,5).Value = "insurance"
$ws.Cells.Item(3,5).Value = "insurance"
$ws.Cells.Item(4,5).Value = "insurance"
$ws.Cells.Item(5,5).Value = "insurance"

# ---- Column F (category) for data rows 2:5 ----
$ws.Cells.Item(2,6).Value = "normal"
$ws.Cells.Item(3,6).Value = "normal"
$ws.Cells.Item(4,6).Value = "normal"
$ws.Cells.Item(5,6).Value = "normal"

# ---- Column G (date) for data rows 2:5 ----
# "2012-04-23" looks like a date, so entering it directly would make Excel
# auto-convert the cell to a date serial number. Instead enter it as a
# formula that evaluates to the literal text, then collapse the formula
# down to its resulting (text) value so the stored cell is a plain string.
$ws.Cells.Item(2,7).Formula = '="2012-04-23"'
$ws.Cells.Item(3,7).Formula = '="2012-04-23"'
$ws.Cells.Item(4,7).Formula = '="2012-04-23"'
$ws.Cells.Item(5,7).Formula = '="2012-04-23"'
$ws.Range("G2:G5").Copy() | Out-Null
$ws.Range("G2:G5").PasteSpecial(-4163) | Out-Null

# ---- Column H (legislator_name) for data rows 2:5 ----
$ws.Cells.Item(2,8).Value = "許忠信"
$ws.Cells.Item(3,8).Value = "許忠信"
$ws.Cells.Item(4,8).Value = "許忠信"
$ws.Cells.Item(5,8).Value = "許忠信"

# ---- Column I (legislator_id) for data rows 2:5 (numeric) ----
$ws.Cells.Item(2,9).Value = 1749
$ws.Cells.Item(3,9).Value = 1749
$ws.Cells.Item(4,9).Value = 1749
$ws.Cells.Item(5,9).Value = 1749

# ---- Column J (source_file) for data rows 2:5 ----
$ws.Cells.Item(2,10).Value = "tmp50641"
$ws.Cells.Item(3,10).Value = "tmp50641"
$ws.Cells.Item(4,10).Value = "tmp50641"
$ws.Cells.Item(5,10).Value = "tmp50641"

# ---- Column K (index) for data rows 2:5 (numeric) ----
$ws.Cells.Item(2,11).Value = 109
$ws.Cells.Item(3,11).Value = 110
$ws.Cells.Item(4,11).Value = 111
$ws.Cells.Item(5,11).Value = 112

# copy data-row formatting onto the newly-introduced data cells F2:K5
$ws.Range("E2").Copy() | Out-Null
$ws.Range("F2:K5").PasteSpecial(-4122) | Out-Null
